# Applies the "updates metadata sheet and start processing data" commit:
#  * fills out more examples on the attribute sheet (rows 15-18, and a new
#    row 22 for an "estimated flow" attribute)
#  * fixes a typo in the K-column data validation list ("interget" -> "interger")
#  * removes the now-unused "_lookups" sheet
#  * moves the active tab from "attribute" to "code_definitions"

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsAttr = $wb.Worksheets.Item("attribute")
$wsCode = $wb.Worksheets.Item("code_definitions")

# --- fill in the "storage_type" / "number_type" style of detail for the
#     ratio rows that previously lacked an F-column (measurement scale) ---
$wsAttr.Range("F15").Value = "numeric"
$wsAttr.Range("F16").Value = "numeric"
$wsAttr.Range("F17").Value = "numeric"
$wsAttr.Range("F18").Value = "numeric"

# --- B19 / B20 get the same "NEED HELP HERE" placeholder used elsewhere ---
$wsAttr.Range("B19").Value = "NEED HELP HERE"
$wsAttr.Range("B20").Value = "NEED HELP HERE"

# --- new row 22: "estimated flow in cubic feet per second" attribute ---
$wsAttr.Range("B22").Value = "estimated flow in cubic feet per second"
$wsAttr.Range("C22").Value = "double"
$wsAttr.Range("D22").Value = "ratio"
$wsAttr.Range("F22").Value = "numeric"
$wsAttr.Range("I22").Value = "ratio"
$wsAttr.Range("K22").Value = "interger"
$wsAttr.Range("O22").Value = 0
$wsAttr.Range("P22").Value = 10000000
$wsAttr.Rows.Item(22).RowHeight = 31.5

# --- fix the "interget" typo in the number_type data validation list ---
$wsAttr.Range("K1:K1048576").Validation.Modify(3, 1, 1, '"natural,whole,interger,real"')

# --- remove the now unused "_lookups" sheet ---
$wb.Worksheets.Item("_lookups").Delete()

# --- move selection / active tab: code_definitions becomes active ---
$wsAttr.Range("B6").Select()
$wsCode.Activate()
$wsCode.Range("A18").Select()
